$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing Issue text (append trailing period) ---
$ws.Range("B2").Value = "unrecognisable characters in answers."
$ws.Range("B3").Value = "Player name is displayed at the top of the page when game ends."

# --- 2. Add new column F "Repair Notes" ---
# Clone formatting from column E (border/alignment) down the whole used range first,
# then overwrite the header cell text.
$ws.Range("E1:E17").Copy()
$ws.Range("F1:F17").PasteSpecial(-4122)
$ws.Range("F1").Value = "Repair Notes"
$ws.Columns.Item(6).ColumnWidth = 20.5

# --- 3. Add new bug row (row 4) ---
# Clone formatting from row 3 (bordered content row) down into row 4 first.
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "The game can start even when there are no players added."
$ws.Range("C4").NumberFormat = "mm-dd-yy"
$ws.Range("C4").Value = (Get-Date -Year 2024 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D4").Value = "When the game is started with no players added, no questions are fetched and the finishing leaderboard is displayed."
$ws.Range("E4").Value = "Ongoing"
$ws.Range("F4").Value = ""
$ws.Rows.Item(4).RowHeight = 58

# --- 4. Change header fill color to bright blue ---
$ws.Range("A1:F1").Interior.Color = 15773696

# --- 5. Add a new "WIP" conditional-formatting rule (orange) on column E, highest priority ---
$rng = $ws.Range("E1:E1048576")
$fc = $rng.FormatConditions.Add(9, 7, "WIP")
$fc.Text = "WIP"
$fc.Interior.Color = 49407
$fc.SetFirstPriority()
